$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.691.98'
$ws.Range("E2").Value = '  +2.88%  '
$ws.Range("D3").Value = '2.199.55'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''258.56'
$ws.Range("E5").Value = '  +2.79%  '
$ws.Range("D6").Value = '''83.46'
$ws.Range("E6").Value = '  +11.01%  '
$ws.Range("D7").Value = '''0.614'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.600'
$ws.Range("E9").Value = '  +2.93%  '
$ws.Range("D10").Value = '''44.59'
$ws.Range("E10").Value = '  +10.69%  '
$ws.Range("D11").Value = '''0.0921'
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").Value = '''7.19'
$ws.Range("E12").Value = '  +6.03%  '
$ws.Range("E13").Value = '  +2.64%  '
$ws.Range("D14").Value = '2.528.09'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '''14.37'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '2.177.38'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("D18").Value = '43.604.94'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("E19").Value = '  +1.63%  '
$ws.Range("D20").Value = '''69.74'
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").Value = '''2.36'
$ws.Range("E22").Value = '  +11.71%  '
$ws.Range("D23").Value = '''231.63'
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("D24").Value = '''8.95'
$ws.Range("E24").Value = '  -5.31%  '
$ws.Range("D26").Value = '''10.66'
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("D28").Value = '''39.46'
$ws.Range("E28").Value = '  +4.82%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.26'
$ws.Range("E29").Value = '  +3.01%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.23'
$ws.Range("E30").Value = '  +4.16%  '
$ws.Range("D31").Value = '''174.23'
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").Value = '''20.40'
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("D33").Value = '''0.0859'
$ws.Range("E33").Value = '  +4.65%  '
$ws.Range("D34").Value = '''5.34'
$ws.Range("E34").Value = '  +4.08%  '
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("E36").Value = '  +4.13%  '
$ws.Range("D37").Value = '''4.53'
$ws.Range("E37").Value = '  +7.76%  '
$ws.Range("E38").Value = '  +7.02%  '
$ws.Range("D39").Value = '''12.52'
$ws.Range("E39").Value = '  +4.73%  '
$ws.Range("E40").Value = '  +9.79%  '
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("D42").Value = '''63.08'
$ws.Range("E42").Value = '  +7.43%  '
$ws.Range("E43").Value = '  +6.32%  '
$ws.Range("D44").Value = '''0.199'
$ws.Range("E44").Value = '  +3.30%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''8.32'
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0979'
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").Value = '''99.77'
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("E48").Value = '  +5.80%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '''1.11'
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").Value = '''0.442'
$ws.Range("E50").Value = '  -3.86%  '
$ws.Range("D51").Value = '''1.48'
$ws.Range("E51").Value = '  +8.00%  '

# Re-apply the original (unstyled) number format to the cells where
# a leading apostrophe was used to keep a numeric-looking price as text,
# so no stray text-format style gets attached to those cells.
# D4 is never touched by this edit and keeps the pristine default style,
# so it is a safe formatting donor for the PasteSpecial(xlPasteFormats) below.
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("D51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
